# Update the "(K) Time Slots" sheet: convert the time-slot description
# strings from 12-hour am/pm notation to 24-hour notation, drop the
# per-row vertical-center formatting override, and restore the zoom /
# selection that was captured at save time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("(K) Time Slots")
[void]$ws.Activate()

$timeSlots = @(
    "MW 11:30 - 12:45",
    "MW 13:00 - 14:15",
    "MW 14:30 - 15:45",
    "MW 16:00 - 17:15",
    "MW 17:30 - 18:45",
    "MW 19:00 - 20:15",
    "MW 20:30 - 21:45",
    "WF 11:30 - 12:45",
    "WF 13:00 - 14:15",
    "WF 14:30 - 15:45",
    "WF 16:00 - 17:15",
    "WF 17:30 - 18:45",
    "WF 19:00 - 20:15",
    "WF 20:30 - 21:45",
    "MF 11:30 - 12:45",
    "MF 13:00 - 14:15",
    "MF 14:30 - 15:45",
    "MF 16:00 - 17:15",
    "MF 17:30 - 18:45",
    "MF 19:00 - 20:15",
    "MF 20:30 - 21:45",
    "MWF 07:00 - 07:50",
    "MWF 08:00 - 08:50",
    "MWF 09:00 - 09:50",
    "MWF 10:00 - 10:50",
    "MWF 11:00 - 11:50",
    "MWF 12:00 - 12:50",
    "MWF 13:00 - 13:50",
    "MWF 14:00 - 14:50",
    "MWF 15:00 - 15:50",
    "MWF 16:00 - 16:50",
    "MWF 17:00 - 17:50",
    "MWF 18:00 - 18:50",
    "MWF 19:00 - 19:50",
    "MWF 20:00 - 20:50",
    "MWF 21:00 - 21:50",
    "TR 07:00 - 08:15",
    "TR 08:30 - 09:15",
    "TR 10:00 - 11:15",
    "TR 13:00 - 14:15",
    "TR 14:30 - 15:45",
    "TR 16:00 - 17:15",
    "TR 17:30 - 18:45",
    "TR 19:00 - 20:15",
    "TR 20:30 - 21:45",
    "MW 11:00 - 12:50",
    "MW 13:00 - 14:50",
    "MW 15:00 - 16:50",
    "MW 17:00 - 18:50",
    "MW 18:00 - 19:50",
    "MW 19:00 - 20:50",
    "MW 20:00 - 21:50",
    "WF 11:00 - 12:50",
    "WF 13:00 - 14:50",
    "WF 15:00 - 16:50",
    "WF 17:00 - 18:50",
    "WF 18:00 - 19:50",
    "WF 19:00 - 20:50",
    "WF 20:00 - 21:50",
    "MF 11:00 - 12:50",
    "MF 13:00 - 14:50",
    "MF 15:00 - 16:50",
    "MF 17:00 - 18:50",
    "MF 18:00 - 19:50",
    "MF 19:00 - 20:50",
    "MF 20:00 - 21:50",
    "MWF 07:00 - 08:05",
    "MWF 08:15 - 09:20",
    "MWF 09:30 - 10:35",
    "MWF 10:45 - 11:50",
    "MWF 12:00 - 13:05",
    "MWF 13:15 - 14:20",
    "MWF 14:30 - 15:35",
    "MWF 15:45 - 16:50",
    "MWF 17:00 - 18:05",
    "MWF 18:15 - 19:20",
    "MWF 19:30 - 20:35",
    "MWF 20:45 - 21:50",
    "TR 08:00 - 09:50",
    "TR 10:00 - 11:50",
    "TR 13:00 - 14:50",
    "TR 15:00 - 16:50",
    "TR 17:00 - 18:50",
    "TR 18:00 - 19:50",
    "TR 19:00 - 20:50",
    "TR 20:00 - 21:50"
)

for ($i = 0; $i -lt $timeSlots.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value2 = $timeSlots[$i]
}

# The row-2..87 "Description" cells previously carried an explicit
# vertical-center alignment style; the edit drops that override so the
# cells fall back to the sheet's default (unstyled) format.
$ws.Range("B2:B87").Style = "Normal"

# Restore the view state captured on save (zoom level + active selection).
$ws.Application.ActiveWindow.Zoom = 389
[void]$ws.Range("C3").Select()
